$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("D2").Value = 44511
$ws.Range("H2").Value = "Sin especificar"
$ws.Range("J2").Value = 500
$ws.Range("K2").Value = 900
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = 950
$ws.Range("P2").Value = 950

# --- Row 3: unchanged ---

# --- Row 4 ---
$ws.Range("D4").Value = 44476
$ws.Range("J4").Value = 300
$ws.Range("K4").Value = 1100
$ws.Range("L4").Value = 1200
$ws.Range("M4").Value = 1150
$ws.Range("P4").Value = 1150

# --- Row 5 ---
$ws.Range("D5").Value = 44512

# --- Row 6 ---
$ws.Range("D6").Value = 44460
$ws.Range("H6").Value = "Verde"
$ws.Range("J6").Value = 120
$ws.Range("K6").Value = 2200
$ws.Range("L6").Value = 2300
$ws.Range("M6").Value = 2250
$ws.Range("P6").Value = 2250

# --- Row 7 ---
$ws.Range("D7").Value = 44516

# --- Row 8 ---
$ws.Range("D8").Value = 44505
$ws.Range("J8").Value = 440

# --- Row 9: unchanged ---

# --- Row 10 ---
$ws.Range("D10").Value = 44510
$ws.Range("J10").Value = 600

# --- Row 11 ---
$ws.Range("D11").Value = 44517
$ws.Range("J11").Value = 500
$ws.Range("K11").Value = 800
$ws.Range("L11").Value = 900
$ws.Range("M11").Value = 850
$ws.Range("P11").Value = 850

# --- Row 12: new row appended ---
$ws.Range("A12").Value = 7
$ws.Range("B12").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C12").Value = "Ñuble"
$ws.Range("D12").Value = 44508
$ws.Range("D12").NumberFormat = $ws.Range("D2").NumberFormat
$ws.Range("E12").Value = 16
$ws.Range("F12").Value = 300000000
$ws.Range("G12").Value = "Espárragos"
$ws.Range("H12").Value = "Sin especificar"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 400
$ws.Range("K12").Value = 900
$ws.Range("L12").Value = 1000
$ws.Range("M12").Value = 950
$ws.Range("N12").Value = "`$/kilo"
$ws.Range("O12").Value = "Provincia de Diguillín"
$ws.Range("P12").Value = 950
$ws.Range("Q12").Value = 1
$ws.Range("R12").Value = "Hortaliza"
